$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 70  # H11
$ws.Cells.Item(11, 9).Value = 70  # I11
$ws.Cells.Item(11, 11).Value = 70  # K11
$ws.Cells.Item(11, 13).Value = 70  # M11
$ws.Cells.Item(32, 8).Value = 3210  # H32
$ws.Cells.Item(32, 10).Value = 3946.25  # J32
$ws.Cells.Item(32, 12).Value = 3946.25  # L32
$ws.Cells.Item(32, 14).Value = -4598.25  # N32
$ws.Cells.Item(33, 8).Value = 1255.421  # H33
$ws.Cells.Item(33, 9).Value = 297.81818  # I33
$ws.Cells.Item(33, 11).Value = 297.81818  # K33
$ws.Cells.Item(33, 13).Value = -68.81817999999998  # M33
$ws.Cells.Item(38, 8).Value = 4215.5  # H38
$ws.Cells.Item(38, 9).Value = 467.41666  # I38
$ws.Cells.Item(38, 10).Value = 7428.143  # J38
$ws.Cells.Item(38, 11).Value = 1402.24998  # K38
$ws.Cells.Item(38, 12).Value = 22284.429  # L38
$ws.Cells.Item(38, 13).Value = -1030.24998  # M38
$ws.Cells.Item(38, 14).Value = -23028.429  # N38
$ws.Cells.Item(40, 8).Value = 2115.3845  # H40
$ws.Cells.Item(40, 9).Value = 1600.2  # I40
$ws.Cells.Item(40, 10).Value = 2437.375  # J40
$ws.Cells.Item(40, 11).Value = 1600.2  # K40
$ws.Cells.Item(40, 12).Value = 2437.375  # L40
$ws.Cells.Item(40, 13).Value = -1425.2  # M40
$ws.Cells.Item(40, 14).Value = -2787.375  # N40
$ws.Cells.Item(43, 8).Value = 5353.364  # H43
$ws.Cells.Item(43, 10).Value = 6165.3335  # J43
$ws.Cells.Item(43, 12).Value = 6165.3335  # L43
$ws.Cells.Item(43, 14).Value = -6303.3335  # N43
$ws.Cells.Item(55, 8).Value = 145  # H55
$ws.Cells.Item(55, 9).Value = 147.5  # I55
$ws.Cells.Item(55, 11).Value = 147.5  # K55
$ws.Cells.Item(55, 13).Value = 66.5  # M55
$ws.Cells.Item(64, 8).Value = 7576.5713  # H64
$ws.Cells.Item(64, 9).Value = 4075.4443  # I64
$ws.Cells.Item(64, 11).Value = 4075.4443  # K64
$ws.Cells.Item(64, 13).Value = -3827.4443  # M64
$ws.Cells.Item(67, 8).Value = 7576.5713  # H67
$ws.Cells.Item(67, 9).Value = 4075.4443  # I67
$ws.Cells.Item(67, 11).Value = 4075.4443  # K67
$ws.Cells.Item(67, 13).Value = -3217.4443  # M67
$ws.Cells.Item(70, 8).Value = 7148942  # H70
$ws.Cells.Item(70, 9).Value = 20003038  # I70
$ws.Cells.Item(70, 11).Value = 60009114  # K70
$ws.Cells.Item(70, 13).Value = -60008844  # M70
$ws.Cells.Item(73, 8).Value = 7148942  # H73
$ws.Cells.Item(73, 9).Value = 20003038  # I73
$ws.Cells.Item(73, 11).Value = 60009114  # K73
$ws.Cells.Item(73, 13).Value = -60008178  # M73
$ws.Cells.Item(76, 8).Value = 4996.6  # H76
$ws.Cells.Item(76, 9).Value = 4993  # I76
$ws.Cells.Item(76, 11).Value = 4993  # K76
$ws.Cells.Item(76, 13).Value = -4678  # M76
$ws.Cells.Item(79, 8).Value = 4996.6  # H79
$ws.Cells.Item(79, 9).Value = 4993  # I79
$ws.Cells.Item(79, 11).Value = 4993  # K79
$ws.Cells.Item(79, 13).Value = -3901  # M79
$ws.Cells.Item(81, 8).Value = 75000  # H81
$ws.Cells.Item(81, 10).Value = 75000  # J81
$ws.Cells.Item(81, 12).Value = 75000  # L81
$ws.Cells.Item(81, 14).Value = -76996  # N81
$ws.Cells.Item(84, 8).Value = 75000  # H84
$ws.Cells.Item(84, 10).Value = 75000  # J84
$ws.Cells.Item(84, 12).Value = 225000  # L84
$ws.Cells.Item(84, 14).Value = -234984  # N84
$ws.Cells.Item(86, 8).Value = 4284.4165  # H86
$ws.Cells.Item(86, 9).Value = 4424.5713  # I86
$ws.Cells.Item(86, 11).Value = 4424.5713  # K86
$ws.Cells.Item(86, 13).Value = -3301.5713  # M86
$ws.Cells.Item(88, 8).Value = 3542.375  # H88
$ws.Cells.Item(88, 9).Value = 609.3333  # I88
$ws.Cells.Item(88, 10).Value = 4219.231  # J88
$ws.Cells.Item(88, 11).Value = 609.3333  # K88
$ws.Cells.Item(88, 12).Value = 4219.231  # L88
$ws.Cells.Item(88, 13).Value = -203.3333  # M88
$ws.Cells.Item(88, 14).Value = -5031.231  # N88
$ws.Cells.Item(89, 8).Value = 4284.4165  # H89
$ws.Cells.Item(89, 9).Value = 4424.5713  # I89
$ws.Cells.Item(89, 11).Value = 22122.8565  # K89
$ws.Cells.Item(89, 13).Value = -16506.8565  # M89
$ws.Cells.Item(91, 8).Value = 3542.375  # H91
$ws.Cells.Item(91, 9).Value = 609.3333  # I91
$ws.Cells.Item(91, 10).Value = 4219.231  # J91
$ws.Cells.Item(91, 11).Value = 609.3333  # K91
$ws.Cells.Item(91, 12).Value = 4219.231  # L91
$ws.Cells.Item(91, 13).Value = 794.6667  # M91
$ws.Cells.Item(91, 14).Value = -7027.231  # N91
$ws.Cells.Item(125, 8).Value = 69954.664  # H125
$ws.Cells.Item(125, 9).Value = 4932  # I125
$ws.Cells.Item(125, 11).Value = 44388  # K125
$ws.Cells.Item(125, 13).Value = -41928  # M125
$ws.Cells.Item(129, 8).Value = 1449.5454  # H129
$ws.Cells.Item(129, 9).Value = 774.36365  # I129
$ws.Cells.Item(129, 11).Value = 2323.09095  # K129
$ws.Cells.Item(129, 13).Value = 2676.90905  # M129
$ws.Cells.Item(132, 8).Value = 1641.0857  # H132
$ws.Cells.Item(132, 9).Value = 1273.6774  # I132
$ws.Cells.Item(132, 11).Value = 3821.0322  # K132
$ws.Cells.Item(132, 13).Value = -1291.0322  # M132
$ws.Cells.Item(137, 8).Value = 12379.8  # H137
$ws.Cells.Item(137, 9).Value = 21660.8  # I137
$ws.Cells.Item(137, 10).Value = 3098.8  # J137
$ws.Cells.Item(137, 11).Value = 64982.39999999999  # K137
$ws.Cells.Item(137, 12).Value = 9296.400000000001  # L137
$ws.Cells.Item(137, 13).Value = -62432.39999999999  # M137
$ws.Cells.Item(137, 14).Value = -14396.4  # N137
$ws.Cells.Item(138, 8).Value = 2420.36  # H138
$ws.Cells.Item(138, 9).Value = 1208.2  # I138
$ws.Cells.Item(138, 11).Value = 3624.6  # K138
$ws.Cells.Item(138, 13).Value = 1515.4  # M138
$ws.Cells.Item(141, 8).Value = 8202.5  # H141
$ws.Cells.Item(141, 9).Value = 7344.8335  # I141
$ws.Cells.Item(141, 10).Value = 8570.071  # J141
$ws.Cells.Item(141, 11).Value = 22034.5005  # K141
$ws.Cells.Item(141, 12).Value = 25710.213  # L141
$ws.Cells.Item(141, 13).Value = -16854.5005  # M141
$ws.Cells.Item(141, 14).Value = -36070.213  # N141

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2857.0715  # H61
$ws.Cells.Item(61, 9).Value = 2807.6155  # I61
$ws.Cells.Item(61, 10).Value = 3500  # J61
$ws.Cells.Item(61, 11).Value = 2807.6155  # K61
$ws.Cells.Item(61, 12).Value = 3500  # L61
$ws.Cells.Item(61, 13).Value = -2595.6155  # M61
$ws.Cells.Item(61, 14).Value = -3924  # N61
$ws.Cells.Item(63, 8).Value = 4030.5  # H63
$ws.Cells.Item(63, 9).Value = 2490.8  # I63
$ws.Cells.Item(63, 10).Value = 5955.125  # J63
$ws.Cells.Item(63, 11).Value = 2490.8  # K63
$ws.Cells.Item(63, 12).Value = 5955.125  # L63
$ws.Cells.Item(63, 13).Value = -1804.8  # M63
$ws.Cells.Item(63, 14).Value = -7327.125  # N63
$ws.Cells.Item(66, 8).Value = 4030.5  # H66
$ws.Cells.Item(66, 9).Value = 2490.8  # I66
$ws.Cells.Item(66, 10).Value = 5955.125  # J66
$ws.Cells.Item(66, 11).Value = 12454  # K66
$ws.Cells.Item(66, 12).Value = 29775.625  # L66
$ws.Cells.Item(66, 13).Value = -9022  # M66
$ws.Cells.Item(66, 14).Value = -36639.625  # N66
$ws.Cells.Item(74, 8).Value = 2179.2307  # H74
$ws.Cells.Item(74, 9).Value = 2235.8333  # I74
$ws.Cells.Item(74, 10).Value = 1500  # J74
$ws.Cells.Item(74, 11).Value = 2235.8333  # K74
$ws.Cells.Item(74, 12).Value = 1500  # L74
$ws.Cells.Item(74, 13).Value = -1361.8333  # M74
$ws.Cells.Item(74, 14).Value = -3248  # N74
$ws.Cells.Item(77, 8).Value = 2179.2307  # H77
$ws.Cells.Item(77, 9).Value = 2235.8333  # I77
$ws.Cells.Item(77, 10).Value = 1500  # J77
$ws.Cells.Item(77, 11).Value = 11179.1665  # K77
$ws.Cells.Item(77, 12).Value = 7500  # L77
$ws.Cells.Item(77, 13).Value = -6811.166499999999  # M77
$ws.Cells.Item(77, 14).Value = -16236  # N77
$ws.Cells.Item(88, 8).Value = 1354.5  # H88
$ws.Cells.Item(88, 9).Value = 610.2222  # I88
$ws.Cells.Item(88, 10).Value = 1963.4546  # J88
$ws.Cells.Item(88, 11).Value = 610.2222  # K88
$ws.Cells.Item(88, 12).Value = 1963.4546  # L88
$ws.Cells.Item(88, 13).Value = -204.2222  # M88
$ws.Cells.Item(88, 14).Value = -2775.4546  # N88
$ws.Cells.Item(91, 8).Value = 1354.5  # H91
$ws.Cells.Item(91, 9).Value = 610.2222  # I91
$ws.Cells.Item(91, 10).Value = 1963.4546  # J91
$ws.Cells.Item(91, 11).Value = 610.2222  # K91
$ws.Cells.Item(91, 12).Value = 1963.4546  # L91
$ws.Cells.Item(91, 13).Value = 793.7778  # M91
$ws.Cells.Item(91, 14).Value = -4771.4546  # N91
$ws.Cells.Item(132, 8).Value = 3445.25  # H132
$ws.Cells.Item(132, 9).Value = 2724.1538  # I132
$ws.Cells.Item(132, 10).Value = 4784.4287  # J132
$ws.Cells.Item(132, 11).Value = 8172.4614  # K132
$ws.Cells.Item(132, 12).Value = 14353.2861  # L132
$ws.Cells.Item(132, 13).Value = -5642.4614  # M132
$ws.Cells.Item(132, 14).Value = -19413.2861  # N132
$ws.Cells.Item(136, 8).Value = 2857.0715  # H136
$ws.Cells.Item(136, 9).Value = 2807.6155  # I136
$ws.Cells.Item(136, 10).Value = 3500  # J136
$ws.Cells.Item(136, 11).Value = 8422.8465  # K136
$ws.Cells.Item(136, 12).Value = 10500  # L136
$ws.Cells.Item(136, 13).Value = -5872.8465  # M136
$ws.Cells.Item(136, 14).Value = -15600  # N136

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 2564.8215  # H94
$ws.Cells.Item(94, 9).Value = 2146.45  # I94
$ws.Cells.Item(94, 10).Value = 3610.75  # J94
$ws.Cells.Item(94, 11).Value = 2146.45  # K94
$ws.Cells.Item(94, 12).Value = 3610.75  # L94
$ws.Cells.Item(94, 13).Value = -1695.45  # M94
$ws.Cells.Item(94, 14).Value = -4512.75  # N94

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(13, 8).Value = 5002.5  # H13
$ws.Cells.Item(13, 9).Value = 10000  # I13
$ws.Cells.Item(13, 10).Value = 5  # J13
$ws.Cells.Item(13, 11).Value = 10000  # K13
$ws.Cells.Item(13, 12).Value = 5  # L13
$ws.Cells.Item(13, 13).Value = -9861  # M13
$ws.Cells.Item(13, 14).Value = -283  # N13
$ws.Cells.Item(31, 8).Value = 5779  # H31
$ws.Cells.Item(31, 9).Value = 1890.9  # I31
$ws.Cells.Item(31, 11).Value = 1890.9  # K31
$ws.Cells.Item(31, 13).Value = -1595.9  # M31
$ws.Cells.Item(34, 8).Value = 5779  # H34
$ws.Cells.Item(34, 9).Value = 1890.9  # I34
$ws.Cells.Item(34, 11).Value = 1890.9  # K34
$ws.Cells.Item(34, 13).Value = -1688.9  # M34
$ws.Cells.Item(58, 8).Value = 1484.7368  # H58
$ws.Cells.Item(58, 9).Value = 1691.4445  # I58
$ws.Cells.Item(58, 10).Value = 1298.7  # J58
$ws.Cells.Item(58, 11).Value = 1691.4445  # K58
$ws.Cells.Item(58, 12).Value = 1298.7  # L58
$ws.Cells.Item(58, 13).Value = -1488.4445  # M58
$ws.Cells.Item(58, 14).Value = -1704.7  # N58
$ws.Cells.Item(86, 8).Value = 0  # H86
$ws.Cells.Item(86, 9).Value = 0  # I86
$ws.Cells.Item(86, 11).Value = 0  # K86
$ws.Cells.Item(86, 13).ClearContents()  # M86
$ws.Cells.Item(89, 8).Value = 0  # H89
$ws.Cells.Item(89, 9).Value = 0  # I89
$ws.Cells.Item(89, 11).Value = 0  # K89
$ws.Cells.Item(89, 13).ClearContents()  # M89
$ws.Cells.Item(99, 8).Value = 2495  # H99
$ws.Cells.Item(99, 9).Value = 2431  # I99
$ws.Cells.Item(99, 11).Value = 2431  # K99
$ws.Cells.Item(99, 13).Value = -933  # M99
$ws.Cells.Item(126, 8).Value = 2495  # H126
$ws.Cells.Item(126, 9).Value = 2431  # I126
$ws.Cells.Item(126, 11).Value = 7293  # K126
$ws.Cells.Item(126, 13).Value = -4823  # M126
$ws.Cells.Item(132, 8).Value = 2054.2  # H132
$ws.Cells.Item(132, 9).Value = 2021.1034  # I132
$ws.Cells.Item(132, 10).Value = 3014  # J132
$ws.Cells.Item(132, 11).Value = 6063.3102  # K132
$ws.Cells.Item(132, 12).Value = 9042  # L132
$ws.Cells.Item(132, 13).Value = -3533.3102  # M132
$ws.Cells.Item(132, 14).Value = -14102  # N132
$ws.Cells.Item(134, 8).Value = 2737.3809  # H134
$ws.Cells.Item(134, 9).Value = 2721.1538  # I134
$ws.Cells.Item(134, 11).Value = 8163.4614  # K134
$ws.Cells.Item(134, 13).Value = -5628.4614  # M134
$ws.Cells.Item(136, 8).Value = 1484.7368  # H136
$ws.Cells.Item(136, 9).Value = 1691.4445  # I136
$ws.Cells.Item(136, 10).Value = 1298.7  # J136
$ws.Cells.Item(136, 11).Value = 5074.333500000001  # K136
$ws.Cells.Item(136, 12).Value = 3896.1  # L136
$ws.Cells.Item(136, 13).Value = -2524.333500000001  # M136
$ws.Cells.Item(136, 14).Value = -8996.1  # N136

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 2080.4375  # H23
$ws.Cells.Item(23, 9).Value = 2350.2222  # I23
$ws.Cells.Item(23, 11).Value = 7050.6666  # K23
$ws.Cells.Item(23, 13).Value = -6815.6666  # M23
$ws.Cells.Item(60, 8).Value = 33433494  # H60
$ws.Cells.Item(60, 9).Value = 37148324  # I60
$ws.Cells.Item(60, 10).Value = 23  # J60
$ws.Cells.Item(60, 11).Value = 111444972  # K60
$ws.Cells.Item(60, 12).Value = 69  # L60
$ws.Cells.Item(60, 13).Value = -111444721  # M60
$ws.Cells.Item(60, 14).Value = -571  # N60
$ws.Cells.Item(113, 8).Value = 1457.1538  # H113
$ws.Cells.Item(113, 10).Value = 1679.3  # J113
$ws.Cells.Item(113, 12).Value = 5037.9  # L113
$ws.Cells.Item(113, 14).Value = -9377.9  # N113
$ws.Cells.Item(136, 8).Value = 2580  # H136
$ws.Cells.Item(136, 9).Value = 1542.8572  # I136
$ws.Cells.Item(136, 11).Value = 4628.571599999999  # K136
$ws.Cells.Item(136, 13).Value = 471.4284000000007  # M136
$ws.Cells.Item(137, 8).Value = 5518.6924  # H137
$ws.Cells.Item(137, 9).Value = 1506.25  # I137
$ws.Cells.Item(137, 10).Value = 11938.6  # J137
$ws.Cells.Item(137, 11).Value = 4518.75  # K137
$ws.Cells.Item(137, 12).Value = 35815.8  # L137
$ws.Cells.Item(137, 13).Value = 581.25  # M137
$ws.Cells.Item(137, 14).Value = -46015.8  # N137
$ws.Cells.Item(140, 8).Value = 60467.117  # H140
$ws.Cells.Item(140, 10).Value = 4225  # J140
$ws.Cells.Item(140, 12).Value = 12675  # L140
$ws.Cells.Item(140, 14).Value = -23035  # N140

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 5744.121  # H122
$ws.Cells.Item(122, 9).Value = 5674.3794  # I122
$ws.Cells.Item(122, 11).Value = 17023.1382  # K122
$ws.Cells.Item(122, 13).Value = -14573.1382  # M122
$ws.Cells.Item(126, 8).Value = 3861.6667  # H126
$ws.Cells.Item(126, 9).Value = 3861.6667  # I126
$ws.Cells.Item(126, 11).Value = 11585.0001  # K126
$ws.Cells.Item(126, 13).Value = -9115.000100000001  # M126
$ws.Cells.Item(132, 8).Value = 2508.739  # H132
$ws.Cells.Item(132, 9).Value = 2621.2104  # I132
$ws.Cells.Item(132, 10).Value = 1974.5  # J132
$ws.Cells.Item(132, 11).Value = 7863.6312  # K132
$ws.Cells.Item(132, 12).Value = 5923.5  # L132
$ws.Cells.Item(132, 13).Value = -5333.6312  # M132
$ws.Cells.Item(132, 14).Value = -10983.5  # N132

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4598.625  # H7
$ws.Cells.Item(7, 9).Value = 4698.4287  # I7
$ws.Cells.Item(7, 11).Value = 4698.4287  # K7
$ws.Cells.Item(7, 13).Value = -4586.4287  # M7
$ws.Cells.Item(40, 8).Value = 6150.9565  # H40
$ws.Cells.Item(40, 9).Value = 4244.231  # I40
$ws.Cells.Item(40, 11).Value = 4244.231  # K40
$ws.Cells.Item(40, 13).Value = -4108.231  # M40
$ws.Cells.Item(46, 8).Value = 2271.2144  # H46
$ws.Cells.Item(46, 9).Value = 933  # I46
$ws.Cells.Item(46, 10).Value = 2636.182  # J46
$ws.Cells.Item(46, 11).Value = 933  # K46
$ws.Cells.Item(46, 12).Value = 2636.182  # L46
$ws.Cells.Item(46, 13).Value = -745  # M46
$ws.Cells.Item(46, 14).Value = -3012.182  # N46
$ws.Cells.Item(68, 8).Value = 12421.667  # H68
$ws.Cells.Item(68, 10).Value = 16131.833  # J68
$ws.Cells.Item(68, 12).Value = 16131.833  # L68
$ws.Cells.Item(68, 14).Value = -17629.833  # N68
$ws.Cells.Item(71, 8).Value = 12421.667  # H71
$ws.Cells.Item(71, 10).Value = 16131.833  # J71
$ws.Cells.Item(71, 12).Value = 80659.16500000001  # L71
$ws.Cells.Item(71, 14).Value = -88147.16500000001  # N71
$ws.Cells.Item(82, 8).Value = 4758.909  # H82
$ws.Cells.Item(82, 9).Value = 4559.6  # I82
$ws.Cells.Item(82, 10).Value = 4925  # J82
$ws.Cells.Item(82, 11).Value = 4559.6  # K82
$ws.Cells.Item(82, 12).Value = 4925  # L82
$ws.Cells.Item(82, 13).Value = -4198.6  # M82
$ws.Cells.Item(82, 14).Value = -5647  # N82
$ws.Cells.Item(85, 8).Value = 4758.909  # H85
$ws.Cells.Item(85, 9).Value = 4559.6  # I85
$ws.Cells.Item(85, 10).Value = 4925  # J85
$ws.Cells.Item(85, 11).Value = 4559.6  # K85
$ws.Cells.Item(85, 12).Value = 4925  # L85
$ws.Cells.Item(85, 13).Value = -3311.6  # M85
$ws.Cells.Item(85, 14).Value = -7421  # N85
$ws.Cells.Item(93, 8).Value = 4326.727  # H93
$ws.Cells.Item(93, 9).Value = 2119  # I93
$ws.Cells.Item(93, 10).Value = 6166.5  # J93
$ws.Cells.Item(93, 11).Value = 2119  # K93
$ws.Cells.Item(93, 12).Value = 6166.5  # L93
$ws.Cells.Item(93, 13).Value = -871  # M93
$ws.Cells.Item(93, 14).Value = -8662.5  # N93
$ws.Cells.Item(95, 8).Value = 6000  # H95
$ws.Cells.Item(95, 10).Value = 6000  # J95
$ws.Cells.Item(95, 12).Value = 6000  # L95
$ws.Cells.Item(95, 14).Value = -11492  # N95
$ws.Cells.Item(100, 8).Value = 5044.385  # H100
$ws.Cells.Item(100, 9).Value = 2739.5715  # I100
$ws.Cells.Item(100, 11).Value = 2739.5715  # K100
$ws.Cells.Item(100, 13).Value = -2198.5715  # M100
$ws.Cells.Item(122, 8).Value = 5730.5  # H122
$ws.Cells.Item(122, 9).Value = 4641.3335  # I122
$ws.Cells.Item(122, 10).Value = 8998  # J122
$ws.Cells.Item(122, 11).Value = 13924.0005  # K122
$ws.Cells.Item(122, 12).Value = 26994  # L122
$ws.Cells.Item(122, 13).Value = -11474.0005  # M122
$ws.Cells.Item(122, 14).Value = -31894  # N122
$ws.Cells.Item(126, 8).Value = 4598.625  # H126
$ws.Cells.Item(126, 9).Value = 4698.4287  # I126
$ws.Cells.Item(126, 11).Value = 14095.2861  # K126
$ws.Cells.Item(126, 13).Value = -11625.2861  # M126
$ws.Cells.Item(136, 8).Value = 5316.4165  # H136
$ws.Cells.Item(136, 9).Value = 5279.7  # I136
$ws.Cells.Item(136, 10).Value = 5500  # J136
$ws.Cells.Item(136, 11).Value = 15839.1  # K136
$ws.Cells.Item(136, 12).Value = 16500  # L136
$ws.Cells.Item(136, 13).Value = -13289.1  # M136
$ws.Cells.Item(136, 14).Value = -21600  # N136

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 10540  # H62
$ws.Cells.Item(62, 10).Value = 12757.143  # J62
$ws.Cells.Item(62, 12).Value = 12757.143  # L62
$ws.Cells.Item(62, 14).Value = -14005.143  # N62
$ws.Cells.Item(65, 8).Value = 10540  # H65
$ws.Cells.Item(65, 10).Value = 12757.143  # J65
$ws.Cells.Item(65, 12).Value = 63785.715  # L65
$ws.Cells.Item(65, 14).Value = -70025.715  # N65
$ws.Cells.Item(107, 8).Value = 3095.0938  # H107
$ws.Cells.Item(107, 9).Value = 1457.963  # I107
$ws.Cells.Item(107, 11).Value = 4373.889  # K107
$ws.Cells.Item(107, 13).Value = -2453.889  # M107
$ws.Cells.Item(113, 8).Value = 1836.8334  # H113
$ws.Cells.Item(113, 10).Value = 1990  # J113
$ws.Cells.Item(113, 12).Value = 5970  # L113
$ws.Cells.Item(113, 14).Value = -10310  # N113
$ws.Cells.Item(122, 8).Value = 6990.857  # H122
$ws.Cells.Item(122, 9).Value = 1578.6  # I122
$ws.Cells.Item(122, 11).Value = 4735.799999999999  # K122
$ws.Cells.Item(122, 13).Value = -2285.799999999999  # M122
$ws.Cells.Item(132, 8).Value = 6322.857  # H132
$ws.Cells.Item(132, 9).Value = 5172.6  # I132
$ws.Cells.Item(132, 11).Value = 15517.8  # K132
$ws.Cells.Item(132, 13).Value = -12987.8  # M132
$ws.Cells.Item(136, 8).Value = 1884.8889  # H136
$ws.Cells.Item(136, 9).Value = 1446.6154  # I136
$ws.Cells.Item(136, 11).Value = 4339.8462  # K136
$ws.Cells.Item(136, 13).Value = -1789.8462  # M136
